$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.913.32'
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").Value = '1.645.57'
$ws.Range("E3").Value = '  +0.57%  '
$ws.Range("E4").Value = '  +0.44%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.39'
$ws.Range("E5").Value = '  -0.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5105'
$ws.Range("E6").Value = '  +1.80%  '
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2574'
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06416'
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.71'
$ws.Range("E10").Value = '  +0.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07776'
$ws.Range("E11").Value = '  +1.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.307'
$ws.Range("E12").Value = '  +1.48%  '
$ws.Range("D13").Value = '1.656.19'
$ws.Range("E13").Value = '  +1.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5473'
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("D15").Value = '0.0₅7892'
$ws.Range("E15").Value = '  -0.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.96'
$ws.Range("E16").Value = '  +2.29%  '
$ws.Range("D17").Value = '25.998.30'
$ws.Range("E17").Value = '  +0.48%  '
$ws.Range("E18").Value = '  +0.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '197.64'
$ws.Range("E19").Value = '  -2.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.424'
$ws.Range("E20").Value = '  +2.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.03'
$ws.Range("E21").Value = '  +0.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.064'
$ws.Range("E22").Value = '  +1.43%  '
$ws.Range("E23").Value = '  +0.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.853'
$ws.Range("E24").Value = '  -4.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '140.71'
$ws.Range("E25").Value = '  -0.40%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1147'
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.891'
$ws.Range("E27").Value = '  +2.65%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.76'
$ws.Range("E28").Value = '  +0.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.240'
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05013'
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.275'
$ws.Range("E31").Value = '  +0.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.202'
$ws.Range("E32").Value = '  +0.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.543'
$ws.Range("E33").Value = '  +0.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.367'
$ws.Range("E34").Value = '  +0.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.8945'
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.588'
$ws.Range("E36").Value = '  -0.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5553'
$ws.Range("E37").Value = '  -1.25%  '
$ws.Range("D38").Value = '1.131.96'
$ws.Range("E38").Value = '  -3.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01566'
$ws.Range("E39").Value = '  +0.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.007'
$ws.Range("E40").Value = '  +0.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.657'
$ws.Range("E41").Value = '  -0.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8151'
$ws.Range("E42").Value = '  +1.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.65'
$ws.Range("E43").Value = '  +0.17%  '
$ws.Range("D44").Value = '0.0₈125'
$ws.Range("E44").Value = '  +8.25%  '
$ws.Range("D45").Value = '1.785.40'
$ws.Range("E45").Value = '  +0.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4536'
$ws.Range("E46").Value = '  +0.43%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.007'
$ws.Range("E47").Value = '  +0.39%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.26'
$ws.Range("E48").Value = '  +0.84%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05093'
$ws.Range("E49").Value = '  +0.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.007'
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.09559'
